$p = $ppt.ActivePresentation

# Remove the "Kubernetes Profile" slide (SlideID 276). The deck is being
# trimmed down to just the "Online Boutique" microservices diagram
# (SlideID 277), so walk the Slides collection and delete the slide whose
# SlideID matches the one being dropped, rather than assuming a fixed
# positional index.
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $s = $p.Slides.Item($i)
    if ($s.SlideID -eq 276) {
        $s.Delete()
    }
}
